# Generate Report for Handoff
# Swaps the "6b361d66..." and "be62f8f9..." rows across the Overview, zh-cn
# and de-de sheets: the 6b361d66 file is now "Ready for handoff" (new
# handoff timestamp), while be62f8f9 remains "Handed back: in sync with
# en-US" but now sorts first.

$wb = $excel.ActiveWorkbook

$sixB  = "6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.md"
$beF   = "be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.md"

$handedBack = "Handed back: in sync with en-US"
$readyForHandoff = "Ready for handoff"

$baseUrl6b = "https://github.com/OpenLocalizationTest/oltest/blob/1f57b83ef396ab8cf88dfabda05ed44ae9884892/e2e/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.md"
$baseUrlBe = "https://github.com/OpenLocalizationTest/oltest/blob/1f57b83ef396ab8cf88dfabda05ed44ae9884892/e2e/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $beF
$ws.Range("B2").Value = $handedBack
$ws.Range("C2").Value = $handedBack
$ws.Range("D2").Value = "2016-03-25 01:00:50"

$ws.Range("A3").Value = $sixB
$ws.Range("B3").Value = $readyForHandoff
$ws.Range("C3").Value = $readyForHandoff
$ws.Range("D3").Value = "2016-03-25 01:02:47"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl6b, "", "", $beF)
$ws.Hyperlinks.Add($ws.Range("A3"), $baseUrlBe, "", "", $sixB)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$handoffUrl6b  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02bd3099f5923ec5e0df8fc76ee18f80b8708c7f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.f4d7b3925ffc053903e7847fcb86813d7f667a32.zh-cn.xlf"
$sourceUrl6b   = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/98b55629da357e46cb2b7200f9e77cd06d22e3b4/e2e/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.md"
$handbackUrl6b = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/45267a5f2e53a26a1276e1d2e9f0efe761f2dc2f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.f4d7b3925ffc053903e7847fcb86813d7f667a32.zh-cn.xlf"

$handoffUrlBe  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02bd3099f5923ec5e0df8fc76ee18f80b8708c7f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.8d65ea536b6db5cde7667a51510fa17ab3ed7fb6.zh-cn.xlf"
$sourceUrlBe   = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/98b55629da357e46cb2b7200f9e77cd06d22e3b4/e2e/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.md"
$handbackUrlBe = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/45267a5f2e53a26a1276e1d2e9f0efe761f2dc2f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.8d65ea536b6db5cde7667a51510fa17ab3ed7fb6.zh-cn.xlf"

$zhHandoffFileBe = "be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.8d65ea536b6db5cde7667a51510fa17ab3ed7fb6.zh-cn.xlf"
$zhHandoffFile6b = "6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.f4d7b3925ffc053903e7847fcb86813d7f667a32.zh-cn.xlf"

# Row 2 now carries the be62f8f9 data (still "Handed back")
$ws.Range("A2").Value = $beF
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $handedBack
$ws.Range("D2").Value = $zhHandoffFileBe
$ws.Range("E2").Value = "2016-03-25 01:00:45"
$ws.Range("F2").Value = $beF
$ws.Range("G2").Value = $zhHandoffFileBe
$ws.Range("H2").Value = "2016-03-25 01:01:34"
$ws.Range("J2").Value = "Include"

# Row 3 now carries the 6b361d66 data, now "Ready for handoff"
$ws.Range("A3").Value = $sixB
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $readyForHandoff
$ws.Range("D3").Value = $zhHandoffFile6b
$ws.Range("E3").Value = "2016-03-25 01:02:43"
$ws.Range("F3").Value = $sixB
$ws.Range("G3").Value = $zhHandoffFile6b
$ws.Range("H3").Value = "2016-03-25 01:01:34"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $sourceUrl6b, "", "", $beF)
$ws.Hyperlinks.Add($ws.Range("D2"), $handoffUrl6b, "", "", $zhHandoffFileBe)
$ws.Hyperlinks.Add($ws.Range("F2"), $baseUrl6b, "", "", $beF)
$ws.Hyperlinks.Add($ws.Range("G2"), $handbackUrl6b, "", "", $zhHandoffFileBe)
$ws.Hyperlinks.Add($ws.Range("A3"), $sourceUrlBe, "", "", $sixB)
$ws.Hyperlinks.Add($ws.Range("D3"), $handoffUrlBe, "", "", $zhHandoffFile6b)
$ws.Hyperlinks.Add($ws.Range("F3"), $baseUrlBe, "", "", $sixB)
$ws.Hyperlinks.Add($ws.Range("G3"), $handbackUrlBe, "", "", $zhHandoffFile6b)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$handoffUrl6bD  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e399cdcde6ef25a0b3ef4172960923f490e4836/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.f4d7b3925ffc053903e7847fcb86813d7f667a32.de-de.xlf"
$sourceUrl6bD   = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ba9bd1265a17d236c0b0d09177683b56d8f7f92b/e2e/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.md"
$handbackUrl6bD = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0b8b03a8caf058d81edd134eb010f08316bfa8bd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.f4d7b3925ffc053903e7847fcb86813d7f667a32.de-de.xlf"

$handoffUrlBeD  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e399cdcde6ef25a0b3ef4172960923f490e4836/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.8d65ea536b6db5cde7667a51510fa17ab3ed7fb6.de-de.xlf"
$sourceUrlBeD   = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ba9bd1265a17d236c0b0d09177683b56d8f7f92b/e2e/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.md"
$handbackUrlBeD = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0b8b03a8caf058d81edd134eb010f08316bfa8bd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.8d65ea536b6db5cde7667a51510fa17ab3ed7fb6.de-de.xlf"

$deHandoffFileBe = "be62f8f9-b98f-487a-8d1b-0d5b8b38b6b9.8d65ea536b6db5cde7667a51510fa17ab3ed7fb6.de-de.xlf"
$deHandoffFile6b = "6b361d66-a5f3-41c4-a856-d3fe8b39ef7d.f4d7b3925ffc053903e7847fcb86813d7f667a32.de-de.xlf"

# Row 2 now carries the be62f8f9 data (still "Handed back")
$ws.Range("A2").Value = $beF
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $handedBack
$ws.Range("D2").Value = $deHandoffFileBe
$ws.Range("E2").Value = "2016-03-25 01:00:50"
$ws.Range("F2").Value = $beF
$ws.Range("G2").Value = $deHandoffFileBe
$ws.Range("H2").Value = "2016-03-25 01:01:50"
$ws.Range("J2").Value = "Include"

# Row 3 now carries the 6b361d66 data, now "Ready for handoff"
$ws.Range("A3").Value = $sixB
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $readyForHandoff
$ws.Range("D3").Value = $deHandoffFile6b
$ws.Range("E3").Value = "2016-03-25 01:02:47"
$ws.Range("F3").Value = $sixB
$ws.Range("G3").Value = $deHandoffFile6b
$ws.Range("H3").Value = "2016-03-25 01:01:50"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $sourceUrl6bD, "", "", $beF)
$ws.Hyperlinks.Add($ws.Range("D2"), $handoffUrl6bD, "", "", $deHandoffFileBe)
$ws.Hyperlinks.Add($ws.Range("F2"), $baseUrl6b, "", "", $beF)
$ws.Hyperlinks.Add($ws.Range("G2"), $handbackUrl6bD, "", "", $deHandoffFileBe)
$ws.Hyperlinks.Add($ws.Range("A3"), $sourceUrlBeD, "", "", $sixB)
$ws.Hyperlinks.Add($ws.Range("D3"), $handoffUrlBeD, "", "", $deHandoffFile6b)
$ws.Hyperlinks.Add($ws.Range("F3"), $baseUrlBe, "", "", $sixB)
$ws.Hyperlinks.Add($ws.Range("G3"), $handbackUrlBeD, "", "", $deHandoffFile6b)
